$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 98 (U17) first, since it is below row 32/93 edits and won't
# disturb their addressing if we do it before the IC8 (row 33) deletion.
$ws.Rows.Item(98).Delete()

# Update row 93 (the TXB0108 DQS_R-PUSON-N20 group row) to absorb U17.
$ws.Range("A93").Value = "U12, U17, U19, U20, U21, U22"
$ws.Range("H93").Value = 6

# Update row 32 (the 74LVC1G3157 mux/demux row) to add IC8 and IC9.
$ws.Range("A32").Value = "IC4, IC5, IC6, IC7, IC8, IC9"
$ws.Range("H32").Value = 6

# Delete row 33 (IC8 inverter row) - the U17_OE inverter has been removed.
$ws.Rows.Item(33).Delete()
